$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Date" column (column E), which also removes its custom
# numeric-format style usage. This shifts Year/Month (old F/G) left to E/F.
$ws.Columns.Item(5).Delete()

# Insert a new "Day" column after Month (now column F), matching the
# Julian day-of-month portion of the old Date value (2014-07-27 -> 27).
$ws.Columns.Item(7).Insert()
$ws.Cells.Item(1, 7).Value = "Day"
$ws.Cells.Item(2, 7).Value = 27

# Remove the trailing taxonomy columns (Genus, Common_Name, Family_ACGC,
# Genus_ACGC, Common_Name_ACGC), which are now columns 24-28.
$ws.Range($ws.Cells.Item(1, 24), $ws.Cells.Item(1, 28)).EntireColumn.Delete()

# Reset selection back to the top-left default cell.
$ws.Range("A1").Select()
